# daily auto push: 2026-02-19 14:13 UTC
# Insert a new data row at row 832 (date 2026/02/19, weekday 木, time 19, ranking 201).
# All subsequent rows shift down by one (old row 873 becomes row 874), and the
# sheet's used dimension grows from A1:D873 to A1:D874.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 832, pushing everything below down by one.
$ws.Rows("832:832").Insert()

# Column A holds a date-like string (e.g. "2026/02/19") that must stay plain text,
# exactly like all the other date cells in the sheet, instead of being auto-converted
# into a date serial number by Excel's input parsing. Temporarily force a text
# number format while assigning the value, then clear the formatting again so the
# cell ends up as an unstyled text cell (matching the rest of the column).
$ws.Range("A832").NumberFormat = "@"
$ws.Range("A832").Value = "2026/02/19"
$ws.Range("A832").ClearFormats()

$ws.Range("B832").Value = "木"
$ws.Range("C832").Value = 19
$ws.Range("D832").Value = 201
